# Commit: Added implementation of MSM measure.
#
# The "interfaceOperations" sheet only listed the operations declared
# directly on each controller interface/class. The MSM (Measure of
# Structural Modification / Member Sharing Measure) computation needs the
# full inherited operation set, so the 9 java.lang.Object operations
# (equals, hashCode, toString, getClass, notify, notifyAll, wait,
# wait(long), wait(long,int)) are now listed for every interface as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("interfaceOperations")

$ws.Range("A2").Value = "org.andante.forum.controller.TopicController"
$ws.Range("B2").Value = "equals(java.lang.Object)"
$ws.Range("C2").Value = "public"
$ws.Range("D2").Value = "boolean"

$ws.Range("A3").Value = "org.andante.forum.controller.TopicController"
$ws.Range("B3").Value = "getTopic(java.lang.Long)"
$ws.Range("C3").Value = "public"
$ws.Range("D3").Value = "org.springframework.http.ResponseEntity"

$ws.Range("A4").Value = "org.andante.forum.controller.TopicController"
$ws.Range("B4").Value = "getHierarchy(java.lang.Long)"
$ws.Range("C4").Value = "public"
$ws.Range("D4").Value = "org.springframework.http.ResponseEntity"

$ws.Range("A5").Value = "org.andante.forum.controller.TopicController"
$ws.Range("B5").Value = "modify(dto.topic.TopicInputDTO)"
$ws.Range("C5").Value = "public"
$ws.Range("D5").Value = "org.springframework.http.ResponseEntity"

$ws.Range("A6").Value = "org.andante.forum.controller.TopicController"
$ws.Range("B6").Value = "toString()"
$ws.Range("C6").Value = "public"
$ws.Range("D6").Value = "java.lang.String"

$ws.Range("A7").Value = "org.andante.forum.controller.TopicController"
$ws.Range("B7").Value = "getSubtopics(java.lang.Long)"
$ws.Range("C7").Value = "public"
$ws.Range("D7").Value = "org.springframework.http.ResponseEntity"

$ws.Range("A8").Value = "org.andante.forum.controller.TopicController"
$ws.Range("B8").Value = "getClass()"
$ws.Range("C8").Value = "public"
$ws.Range("D8").Value = "java.lang.Class"

$ws.Range("A9").Value = "org.andante.forum.controller.TopicController"
$ws.Range("B9").Value = "getTop(java.lang.Integer, java.lang.Integer)"
$ws.Range("C9").Value = "public"
$ws.Range("D9").Value = "org.springframework.http.ResponseEntity"

$ws.Range("A10").Value = "org.andante.forum.controller.TopicController"
$ws.Range("B10").Value = "notifyAll()"
$ws.Range("C10").Value = "public"
$ws.Range("D10").Value = "void"

$ws.Range("A11").Value = "org.andante.forum.controller.TopicController"
$ws.Range("B11").Value = "hashCode()"
$ws.Range("C11").Value = "public"
$ws.Range("D11").Value = "int"

$ws.Range("A12").Value = "org.andante.forum.controller.TopicController"
$ws.Range("B12").Value = "wait()"
$ws.Range("C12").Value = "public"
$ws.Range("D12").Value = "void"

$ws.Range("A13").Value = "org.andante.forum.controller.TopicController"
$ws.Range("B13").Value = "TopicController(org.andante.forum.logic.service.impl.TopicServiceImpl, org.andante.forum.controller.mapper.TopicDTOModelMapper, org.andante.mappers.OperationHttpStatusMapper)"
$ws.Range("C13").Value = "public"
$ws.Range("D13").Value = "void"

$ws.Range("A14").Value = "org.andante.forum.controller.TopicController"
$ws.Range("B14").Value = "getPage(dto.topic.TopicQuerySpecification)"
$ws.Range("C14").Value = "public"
$ws.Range("D14").Value = "org.springframework.http.ResponseEntity"

$ws.Range("A15").Value = "org.andante.forum.controller.TopicController"
$ws.Range("B15").Value = "create(dto.topic.TopicInputDTO)"
$ws.Range("C15").Value = "public"
$ws.Range("D15").Value = "org.springframework.http.ResponseEntity"

$ws.Range("A16").Value = "org.andante.forum.controller.TopicController"
$ws.Range("B16").Value = "notify()"
$ws.Range("C16").Value = "public"
$ws.Range("D16").Value = "void"

$ws.Range("A17").Value = "org.andante.forum.controller.TopicController"
$ws.Range("B17").Value = "wait(long)"
$ws.Range("C17").Value = "public"
$ws.Range("D17").Value = "void"

$ws.Range("A18").Value = "org.andante.forum.controller.TopicController"
$ws.Range("B18").Value = "wait(long, int)"
$ws.Range("C18").Value = "public"
$ws.Range("D18").Value = "void"

$ws.Range("A19").Value = "org.andante.forum.controller.TopicController"
$ws.Range("B19").Value = "delete(java.lang.Long)"
$ws.Range("C19").Value = "public"
$ws.Range("D19").Value = "org.springframework.http.ResponseEntity"

$ws.Range("A20").Value = "org.andante.forum.controller.PostController"
$ws.Range("B20").Value = "getPage(dto.post.PostQuerySpecification)"
$ws.Range("C20").Value = "public"
$ws.Range("D20").Value = "org.springframework.http.ResponseEntity"

$ws.Range("A21").Value = "org.andante.forum.controller.PostController"
$ws.Range("B21").Value = "like(dto.post.PostLikeDTO)"
$ws.Range("C21").Value = "public"
$ws.Range("D21").Value = "org.springframework.http.ResponseEntity"

$ws.Range("A22").Value = "org.andante.forum.controller.PostController"
$ws.Range("B22").Value = "equals(java.lang.Object)"
$ws.Range("C22").Value = "public"
$ws.Range("D22").Value = "boolean"

$ws.Range("A23").Value = "org.andante.forum.controller.PostController"
$ws.Range("B23").Value = "toString()"
$ws.Range("C23").Value = "public"
$ws.Range("D23").Value = "java.lang.String"

$ws.Range("A24").Value = "org.andante.forum.controller.PostController"
$ws.Range("B24").Value = "getTopPage(dto.post.TopQuerySpecification)"
$ws.Range("C24").Value = "public"
$ws.Range("D24").Value = "org.springframework.http.ResponseEntity"

$ws.Range("A25").Value = "org.andante.forum.controller.PostController"
$ws.Range("B25").Value = "getClass()"
$ws.Range("C25").Value = "public"
$ws.Range("D25").Value = "java.lang.Class"

$ws.Range("A26").Value = "org.andante.forum.controller.PostController"
$ws.Range("B26").Value = "notifyAll()"
$ws.Range("C26").Value = "public"
$ws.Range("D26").Value = "void"

$ws.Range("A27").Value = "org.andante.forum.controller.PostController"
$ws.Range("B27").Value = "modify(dto.post.PostInputDTO)"
$ws.Range("C27").Value = "public"
$ws.Range("D27").Value = "org.springframework.http.ResponseEntity"

$ws.Range("A28").Value = "org.andante.forum.controller.PostController"
$ws.Range("B28").Value = "hashCode()"
$ws.Range("C28").Value = "public"
$ws.Range("D28").Value = "int"

$ws.Range("A29").Value = "org.andante.forum.controller.PostController"
$ws.Range("B29").Value = "wait()"
$ws.Range("C29").Value = "public"
$ws.Range("D29").Value = "void"

$ws.Range("A30").Value = "org.andante.forum.controller.PostController"
$ws.Range("B30").Value = "create(dto.post.PostInputDTO)"
$ws.Range("C30").Value = "public"
$ws.Range("D30").Value = "org.springframework.http.ResponseEntity"

$ws.Range("A31").Value = "org.andante.forum.controller.PostController"
$ws.Range("B31").Value = "notify()"
$ws.Range("C31").Value = "public"
$ws.Range("D31").Value = "void"

$ws.Range("A32").Value = "org.andante.forum.controller.PostController"
$ws.Range("B32").Value = "get(java.lang.Long)"
$ws.Range("C32").Value = "public"
$ws.Range("D32").Value = "org.springframework.http.ResponseEntity"

$ws.Range("A33").Value = "org.andante.forum.controller.PostController"
$ws.Range("B33").Value = "getLikedByUser(java.lang.String)"
$ws.Range("C33").Value = "public"
$ws.Range("D33").Value = "org.springframework.http.ResponseEntity"

$ws.Range("A34").Value = "org.andante.forum.controller.PostController"
$ws.Range("B34").Value = "wait(long)"
$ws.Range("C34").Value = "public"
$ws.Range("D34").Value = "void"

$ws.Range("A35").Value = "org.andante.forum.controller.PostController"
$ws.Range("B35").Value = "PostController(org.andante.forum.logic.service.impl.PostServiceImpl, org.andante.forum.controller.mapper.PostDTOModelMapper, org.andante.forum.controller.mapper.PostLikesDTOModelMapper, org.andante.mappers.OperationHttpStatusMapper)"
$ws.Range("C35").Value = "public"
$ws.Range("D35").Value = "void"

$ws.Range("A36").Value = "org.andante.forum.controller.PostController"
$ws.Range("B36").Value = "wait(long, int)"
$ws.Range("C36").Value = "public"
$ws.Range("D36").Value = "void"

$ws.Range("A37").Value = "org.andante.forum.controller.PostController"
$ws.Range("B37").Value = "delete(java.lang.Long)"
$ws.Range("C37").Value = "public"
$ws.Range("D37").Value = "org.springframework.http.ResponseEntity"

$ws.Range("A38").Value = "org.andante.forum.controller.PostResponseController"
$ws.Range("B38").Value = "like(dto.response.PostResponseLikeDTO)"
$ws.Range("C38").Value = "public"
$ws.Range("D38").Value = "org.springframework.http.ResponseEntity"

$ws.Range("A39").Value = "org.andante.forum.controller.PostResponseController"
$ws.Range("B39").Value = "equals(java.lang.Object)"
$ws.Range("C39").Value = "public"
$ws.Range("D39").Value = "boolean"

$ws.Range("A40").Value = "org.andante.forum.controller.PostResponseController"
$ws.Range("B40").Value = "modify(dto.response.PostResponseInputDTO)"
$ws.Range("C40").Value = "public"
$ws.Range("D40").Value = "org.springframework.http.ResponseEntity"

$ws.Range("A41").Value = "org.andante.forum.controller.PostResponseController"
$ws.Range("B41").Value = "toString()"
$ws.Range("C41").Value = "public"
$ws.Range("D41").Value = "java.lang.String"

$ws.Range("A42").Value = "org.andante.forum.controller.PostResponseController"
$ws.Range("B42").Value = "PostResponseController(org.andante.forum.logic.service.impl.PostResponseServiceImpl, org.andante.forum.controller.mapper.PostResponseDTOModelMapper, org.andante.forum.controller.mapper.PostResponsesLikesDTOModelMapper, org.andante.mappers.OperationHttpStatusMapper)"
$ws.Range("C42").Value = "public"
$ws.Range("D42").Value = "void"

$ws.Range("A43").Value = "org.andante.forum.controller.PostResponseController"
$ws.Range("B43").Value = "getPage(dto.response.PostResponseQuerySpecification)"
$ws.Range("C43").Value = "public"
$ws.Range("D43").Value = "org.springframework.http.ResponseEntity"

$ws.Range("A44").Value = "org.andante.forum.controller.PostResponseController"
$ws.Range("B44").Value = "getClass()"
$ws.Range("C44").Value = "public"
$ws.Range("D44").Value = "java.lang.Class"

$ws.Range("A45").Value = "org.andante.forum.controller.PostResponseController"
$ws.Range("B45").Value = "create(dto.response.PostResponseInputDTO)"
$ws.Range("C45").Value = "public"
$ws.Range("D45").Value = "org.springframework.http.ResponseEntity"

$ws.Range("A46").Value = "org.andante.forum.controller.PostResponseController"
$ws.Range("B46").Value = "notifyAll()"
$ws.Range("C46").Value = "public"
$ws.Range("D46").Value = "void"

$ws.Range("A47").Value = "org.andante.forum.controller.PostResponseController"
$ws.Range("B47").Value = "hashCode()"
$ws.Range("C47").Value = "public"
$ws.Range("D47").Value = "int"

$ws.Range("A48").Value = "org.andante.forum.controller.PostResponseController"
$ws.Range("B48").Value = "wait()"
$ws.Range("C48").Value = "public"
$ws.Range("D48").Value = "void"

$ws.Range("A49").Value = "org.andante.forum.controller.PostResponseController"
$ws.Range("B49").Value = "notify()"
$ws.Range("C49").Value = "public"
$ws.Range("D49").Value = "void"

$ws.Range("A50").Value = "org.andante.forum.controller.PostResponseController"
$ws.Range("B50").Value = "get(java.lang.Long)"
$ws.Range("C50").Value = "public"
$ws.Range("D50").Value = "org.springframework.http.ResponseEntity"

$ws.Range("A51").Value = "org.andante.forum.controller.PostResponseController"
$ws.Range("B51").Value = "wait(long)"
$ws.Range("C51").Value = "public"
$ws.Range("D51").Value = "void"

$ws.Range("A52").Value = "org.andante.forum.controller.PostResponseController"
$ws.Range("B52").Value = "wait(long, int)"
$ws.Range("C52").Value = "public"
$ws.Range("D52").Value = "void"

$ws.Range("A53").Value = "org.andante.forum.controller.PostResponseController"
$ws.Range("B53").Value = "delete(java.lang.Long)"
$ws.Range("C53").Value = "public"
$ws.Range("D53").Value = "org.springframework.http.ResponseEntity"

$ws.Range("A54").Value = "org.andante.forum.controller.UserController"
$ws.Range("B54").Value = "equals(java.lang.Object)"
$ws.Range("C54").Value = "public"
$ws.Range("D54").Value = "boolean"

$ws.Range("A55").Value = "org.andante.forum.controller.UserController"
$ws.Range("B55").Value = "create(dto.UserDTO)"
$ws.Range("C55").Value = "public"
$ws.Range("D55").Value = "org.springframework.http.ResponseEntity"

$ws.Range("A56").Value = "org.andante.forum.controller.UserController"
$ws.Range("B56").Value = "UserController(org.andante.forum.logic.service.impl.UserServiceImpl, org.andante.forum.controller.mapper.UserDTOModelMapper, org.andante.mappers.OperationHttpStatusMapper)"
$ws.Range("C56").Value = "public"
$ws.Range("D56").Value = "void"

$ws.Range("A57").Value = "org.andante.forum.controller.UserController"
$ws.Range("B57").Value = "toString()"
$ws.Range("C57").Value = "public"
$ws.Range("D57").Value = "java.lang.String"

$ws.Range("A58").Value = "org.andante.forum.controller.UserController"
$ws.Range("B58").Value = "delete(java.lang.String)"
$ws.Range("C58").Value = "public"
$ws.Range("D58").Value = "org.springframework.http.ResponseEntity"

$ws.Range("A59").Value = "org.andante.forum.controller.UserController"
$ws.Range("B59").Value = "getClass()"
$ws.Range("C59").Value = "public"
$ws.Range("D59").Value = "java.lang.Class"

$ws.Range("A60").Value = "org.andante.forum.controller.UserController"
$ws.Range("B60").Value = "notifyAll()"
$ws.Range("C60").Value = "public"
$ws.Range("D60").Value = "void"

$ws.Range("A61").Value = "org.andante.forum.controller.UserController"
$ws.Range("B61").Value = "hashCode()"
$ws.Range("C61").Value = "public"
$ws.Range("D61").Value = "int"

$ws.Range("A62").Value = "org.andante.forum.controller.UserController"
$ws.Range("B62").Value = "wait()"
$ws.Range("C62").Value = "public"
$ws.Range("D62").Value = "void"

$ws.Range("A63").Value = "org.andante.forum.controller.UserController"
$ws.Range("B63").Value = "notify()"
$ws.Range("C63").Value = "public"
$ws.Range("D63").Value = "void"

$ws.Range("A64").Value = "org.andante.forum.controller.UserController"
$ws.Range("B64").Value = "wait(long)"
$ws.Range("C64").Value = "public"
$ws.Range("D64").Value = "void"

$ws.Range("A65").Value = "org.andante.forum.controller.UserController"
$ws.Range("B65").Value = "get(java.lang.String)"
$ws.Range("C65").Value = "public"
$ws.Range("D65").Value = "org.springframework.http.ResponseEntity"

$ws.Range("A66").Value = "org.andante.forum.controller.UserController"
$ws.Range("B66").Value = "wait(long, int)"
$ws.Range("C66").Value = "public"
$ws.Range("D66").Value = "void"

